# 4th excel utils pkg added
#
# Renames the existing sheet to "correct log" and adds a new "incorrect
# log" sheet (for the locked_out_user login attempt), mirroring the
# username/password layout of the first sheet.

$wb = $excel.ActiveWorkbook

# Rename Sheet1 -> "correct log"
$ws1 = $wb.ActiveSheet
$ws1.Name = "correct log"

# It is no longer the active tab; update its selection.
$ws1.Range("A1:B1").Select() | Out-Null

# Add the new "incorrect log" sheet right after "correct log".
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "incorrect log"

# Header + locked-out-user login row.
$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"
$ws2.Range("A2").Value = "locked_out_user"
$ws2.Range("B2").Value = "secret_sauce"

# Best-fit the columns to their content, like the first sheet.
$ws2.Columns.Item(1).AutoFit()
$ws2.Columns.Item(2).AutoFit()

$ws2.PageSetup.Orientation = 1

# New sheet becomes the active tab, selection on F6.
$ws2.Range("F6").Select() | Out-Null
$ws2.Activate() | Out-Null
